# Apply updated cryptocurrency price/volume data to sheet1 (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell 'D2' '64.009.08'
Set-TextCell 'E2' '  -3.21%  '
Set-TextCell 'D3' '3.194.76'
Set-TextCell 'E3' '  -2.99%  '
Set-TextCell 'D5' '570.82'
Set-TextCell 'E5' '  -2.56%  '
Set-TextCell 'D6' '170.36'
Set-TextCell 'E6' '  -5.56%  '
Set-TextCell 'D7' '0.616'
Set-TextCell 'E7' '  -6.09%  '
Set-TextCell 'E8' '  +0.04%  '
Set-TextCell 'D9' '3.194.25'
Set-TextCell 'E9' '  -2.97%  '
Set-TextCell 'D10' '0.121'
Set-TextCell 'E10' '  -3.34%  '
Set-TextCell 'D11' '6.75'
Set-TextCell 'E11' '  +0.13%  '
Set-TextCell 'D12' '0.387'
Set-TextCell 'E12' '  -4.56%  '
Set-TextCell 'D13' '3.750.37'
Set-TextCell 'E13' '  -3.03%  '
Set-TextCell 'E14' '  -1.68%  '
Set-TextCell 'D15' '64.190.32'
Set-TextCell 'E15' '  -3.00%  '
Set-TextCell 'D16' '25.57'
Set-TextCell 'E16' '  -3.36%  '
Set-TextCell 'D17' '0.0000159'
Set-TextCell 'E17' '  -2.43%  '
Set-TextCell 'D18' '3.193.29'
Set-TextCell 'E18' '  -2.79%  '
Set-TextCell 'D19' '415.39'
Set-TextCell 'E19' '  -4.75%  '
Set-TextCell 'D20' '5.38'
Set-TextCell 'E20' '  -2.00%  '
Set-TextCell 'D21' '12.85'
Set-TextCell 'E21' '  -3.12%  '
Set-TextCell 'D22' '7.19'
Set-TextCell 'E22' '  -3.32%  '
Set-TextCell 'D23' '0.999'
Set-TextCell 'E23' '  -0.10%  '
Set-TextCell 'D24' '70.94'
Set-TextCell 'E24' '  -2.02%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 'D25' '0.203'
Set-TextCell 'E25' '  +2.59%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 'D26' '0.493'
Set-TextCell 'E26' '  -3.46%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell 'D27' '0.0000110'
Set-TextCell 'E27' '  -2.94%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D28' '8.81'
Set-TextCell 'E28' '  -0.53%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell 'D29' '1.00'
Set-TextCell 'E29' '  -0.15%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D30' '1.85'
Set-TextCell 'E30' '  -5.97%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 'D31' '21.93'
Set-TextCell 'E31' '  -1.81%  '
$ws.Range('B32').Value = 'USDe'
$ws.Range('C32').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell 'D32' '1.00'
Set-TextCell 'E32' '  +0.15%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D33' '5.02'
Set-TextCell 'E33' '  -3.66%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D34' '6.41'
Set-TextCell 'E34' '  -3.04%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell 'D35' '1.14'
Set-TextCell 'E35' '  -4.39%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D36' '156.21'
Set-TextCell 'E36' '  -1.13%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D37' '1.38'
Set-TextCell 'E37' '  -2.65%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 'D38' '2.741.31'
Set-TextCell 'E38' '  -2.06%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell 'D39' '1.71'
Set-TextCell 'E39' '  -3.58%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D40' '25.26'
Set-TextCell 'E40' '  -4.69%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D41' '4.20'
Set-TextCell 'E41' '  -3.50%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D42' '0.720'
Set-TextCell 'E42' '  -6.94%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 'D43' '38.95'
Set-TextCell 'E43' '  -3.20%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D44' '5.78'
Set-TextCell 'E44' '  -4.90%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D45' '0.0629'
Set-TextCell 'E45' '  -4.88%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 'D46' '22.22'
Set-TextCell 'E46' '  -4.44%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D47' '298.51'
Set-TextCell 'E47' '  -6.93%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell 'D48' '2.13'
Set-TextCell 'E48' '  -8.04%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D49' '0.0262'
Set-TextCell 'E49' '  -2.10%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D50' '0.0997'
Set-TextCell 'E50' '  -6.16%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell 'D51' '1.00'
Set-TextCell 'E51' '  +0.03%  '
